# Update countries & provincias Spain
# Refresh the COVID-19 country data table (sheet "Pais") to the
# "12 de Agosto de 2020 a las 13:20" snapshot:
#  - update case counters for a number of countries
#  - re-sort pushed a few countries past their neighbour, so swap the
#    two rows' country name + counters where that happened
#  - bump the "Datos actualizados ..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Row, $Country, $Values) {
    if ($Country) {
        $ws.Cells.Item($Row, 1).Value2 = $Country
    }
    foreach ($col in $Values.Keys) {
        $colIndex = switch ($col) {
            'B' { 2 }
            'C' { 3 }
            'D' { 4 }
            'E' { 5 }
            'F' { 6 }
            'G' { 7 }
            'H' { 8 }
        }
        $ws.Cells.Item($Row, $colIndex).Value2 = $Values[$col]
    }
}

# Rows whose country stays the same, only the counters move.
Set-RowValues 6   $null @{ B = 2333166; C = 4761; D = 1640463; E = 646487 }
Set-RowValues 43  $null @{ B = 65177;   C = 1415; D = 31048;   E = 31322; G = 43; H = 2807 }
Set-RowValues 44  $null @{ B = 63212;   C = 246;  D = 57193;   E = 5661 }
Set-RowValues 58  $null @{ B = 37169;   C = 274;  E = 2478;    G = 1;  H = 1991 }
Set-RowValues 85  $null @{ B = 11587;   C = 207;  D = 7523;    E = 3822; G = 4; H = 242 }
Set-RowValues 95  $null @{ D = 7050;    E = 259 }
Set-RowValues 159 $null @{ B = 880;     C = 14;   E = 463 }
Set-RowValues 174 $null @{ B = 334;     C = 16;   E = 109 }

# Rows where the refreshed counters moved the country past its neighbour,
# so the two rows swap places (country name + all counters).
Set-RowValues 69  'Nepal'          @{ B = 24432; C = 484; D = 16728; E = 7613; G = 8; H = 91 }
Set-RowValues 70  'Etiopia'        @{ B = 24175; C = 0;   D = 10696; E = 13039; G = 0; H = 440 }

Set-RowValues 142 'Uganda'         @{ B = 1332; C = 19; D = 1139; E = 184; H = 9 }
Set-RowValues 143 'Siria'          @{ B = 1327; C = 0;  D = 385;  E = 889; H = 53 }

Set-RowValues 150 'Malta'          @{ B = 1190; C = 49; D = 695;  E = 486; H = 9 }
Set-RowValues 151 'Niger'          @{ B = 1158; C = 0;  D = 1065; E = 24;  H = 69 }

Set-RowValues 213 'Islas Malvinas' @{ D = 13; H = 0 }
Set-RowValues 214 'Montserrat'     @{ D = 12; H = 1 }

# Timestamp footer (row 1)
$ws.Range("A1").Value2 = "Datos actualizados a 12 de Agosto de 2020 a las 13:20"
